$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell (G1, "sum") onto the new
# header cell H1, then overwrite its value with "Save".
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# New "Save" column values for data rows 2-10
$values = @(1, 0, 0, 1, 1, 1, 1, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
